$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64; this shifts existing rows 64-67 down to 65-68
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly data
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value = 44706
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112012
$ws.Cells.Item(64, 7).Value = "Espinaca"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 250
$ws.Cells.Item(64, 11).Value = 2500
$ws.Cells.Item(64, 12).Value = 3000
$ws.Cells.Item(64, 13).Value = 2750
$ws.Cells.Item(64, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 917
$ws.Cells.Item(64, 17).Value = 3
$ws.Cells.Item(64, 18).Value = "Hortaliza"

# Row 66 (formerly row 65's successor becomes row 66 after shift): quality changed Segunda -> Primera,
# and price columns updated. (Row numbers below reflect post-insert positions.)
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 11).Value = 1400
$ws.Cells.Item(66, 12).Value = 1500
$ws.Cells.Item(66, 13).Value = 1450
$ws.Cells.Item(66, 16).Value = 483
